$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 13; $row++) {
    for ($col = 2; $col -le 8; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $excel.WorksheetFunction.Round([double]$cell.Value2, 0)
    }
}
